$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F4").Value = "Fluor Corporate"
$ws.Range("F17").Select()
